$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to keep its value as plain text even when the string
    # looks like a valid number (e.g. "589.05", "2.04", "0.999", ...).
    # We briefly apply a text format, write the value, then restore the
    # cell's original style so no visible formatting changes.
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.451.48"
$ws.Range("E2").Value = "  +2.70%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.485.76"
$ws.Range("E3").Value = "  +1.86%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "589.05"
$ws.Range("E5").Value = "  +2.35%  "

# Row 6 - Solana
Set-TextValue "D6" "167.54"
$ws.Range("E6").Value = "  +1.49%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.483.05"
$ws.Range("E8").Value = "  +1.79%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +6.29%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +0.35%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +5.71%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.32%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.091.21"
$ws.Range("E13").Value = "  +2.33%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.37%  "

# Row 15 - Avalanche
Set-TextValue "D15" "27.86"
$ws.Range("E15").Value = "  +3.10%  "

# Row 16 - was WrappedBTC, now ShibaInu
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000178"
$ws.Range("E16").Value = "  +2.45%  "

# Row 17 - was ShibaInu, now WrappedBTC
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "66.496.44"
$ws.Range("E17").Value = "  +2.76%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.491.72"
$ws.Range("E18").Value = "  +2.89%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.00%  "

# Row 20 - Chainlink
Set-TextValue "D20" "13.90"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "389.90"
$ws.Range("E21").Value = "  +2.74%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.07%  "

# Row 23 - Litecoin
Set-TextValue "D23" "72.66"
$ws.Range("E23").Value = "  +2.46%  "

# Row 24 - Dai
Set-TextValue "D24" "0.999"
$ws.Range("E24").Value = "  -0.28%  "

# Row 25 - Polygon
$ws.Range("E25").Value = "  +2.37%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +3.93%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "10.12"
$ws.Range("E27").Value = "  +3.84%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.07%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.25%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  +2.83%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +1.99%  "

# Row 32 - PancakeSwap
Set-TextValue "D32" "2.04"
$ws.Range("E32").Value = "  +1.65%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "23.58"
$ws.Range("E33").Value = "  +2.67%  "

# Row 34 - Aptos
$ws.Range("E34").Value = "  +3.19%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +6.05%  "

# Row 36 - Monero
Set-TextValue "D36" "162.70"
$ws.Range("E36").Value = "  +1.96%  "

# Row 37 - Mantle
Set-TextValue "D37" "0.893"
$ws.Range("E37").Value = "  +2.69%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +3.42%  "

# Row 39 - RenderToken
Set-TextValue "D39" "6.79"
$ws.Range("E39").Value = "  +4.64%  "

# Row 40 - was Hedera, now Filecoin
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D40" "4.61"
$ws.Range("E40").Value = "  +5.07%  "

# Row 41 - was Filecoin, now Hedera
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.0734"
$ws.Range("E41").Value = "  +1.28%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "26.22"
$ws.Range("E42").Value = "  +1.33%  "

# Row 43 - Maker
Set-TextValue "D43" "2.766.91"
$ws.Range("E43").Value = "  -0.21%  "

# Row 44 - was InjectiveProtocol, now OKB
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D44" "42.74"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45 - was OKB, now InjectiveProtocol
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "26.38"
$ws.Range("E45").Value = "  +2.13%  "

# Row 46 - dogwifhat
$ws.Range("E46").Value = "  +2.66%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +0.97%  "

# Row 48 - Bittensor
Set-TextValue "D48" "340.37"
$ws.Range("E48").Value = "  +3.05%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  +2.19%  "

# Row 50 - Arweave
Set-TextValue "D50" "33.29"
$ws.Range("E50").Value = "  +8.71%  "

# Row 51 - SuiNetwork
Set-TextValue "D51" "0.849"
$ws.Range("E51").Value = "  +3.89%  "
